$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 218, shifting existing rows 218-227 down to 220-229
$ws.Range("A218:R219").EntireRow.Insert()

# Populate new row 218 (Primera)
$ws.Range("A218").Value = 11
$ws.Range("B218").Value = "Vega Monumental Concepción"
$ws.Range("C218").Value = "Bíobío"
$ws.Range("D218").Value = 44509
$ws.Range("E218").Value = 8
$ws.Range("F218").Value = 100112006
$ws.Range("G218").Value = "Repollo"
$ws.Range("H218").Value = "Crespo record"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 2000
$ws.Range("K218").Value = 700
$ws.Range("L218").Value = 800
$ws.Range("M218").Value = 750
$ws.Range("N218").Value = '$/unidad'
$ws.Range("O218").Value = "Región Metropolitana"
$ws.Range("P218").Value = 750
$ws.Range("Q218").Value = 1
$ws.Range("R218").Value = "Hortaliza"

# Populate new row 219 (Segunda)
$ws.Range("A219").Value = 11
$ws.Range("B219").Value = "Vega Monumental Concepción"
$ws.Range("C219").Value = "Bíobío"
$ws.Range("D219").Value = 44509
$ws.Range("E219").Value = 8
$ws.Range("F219").Value = 100112006
$ws.Range("G219").Value = "Repollo"
$ws.Range("H219").Value = "Crespo record"
$ws.Range("I219").Value = "Segunda"
$ws.Range("J219").Value = 1000
$ws.Range("K219").Value = 600
$ws.Range("L219").Value = 600
$ws.Range("M219").Value = 600
$ws.Range("N219").Value = '$/unidad'
$ws.Range("O219").Value = "Región Metropolitana"
$ws.Range("P219").Value = 600
$ws.Range("Q219").Value = 1
$ws.Range("R219").Value = "Hortaliza"
